# GFTA.xlsx update — "Updated data with DX explanations"
#
# The GFTA "part 2" block (sounds-in-sentences) gets its own raw/standard/
# percentile rows (mirroring "part 1"), variable names are renamed to the
# GFTA_SIW_* / GFTA_SIS_* convention, a "Scores" banner row is removed, and
# the trailing thick-bottom-border row is dropped (folded into a plain row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Drop the old row 11 (the heavy "thick bottom border" closing row of
#    the question block). Deleting it (rather than just clearing it)
#    shifts rows 12-18 up to 11-17, which removes the special row height /
#    thick-border row formatting in one shot and lands the sheet on the
#    final 17-row extent.
# ---------------------------------------------------------------------
$ws.Rows.Item(11).Delete()

# ---------------------------------------------------------------------
# 2. Rewrite the variable-name column (B) for the "part 1" (sounds-in-
#    words) question block to the new GFTA_SIW_* naming.
# ---------------------------------------------------------------------
$ws.Range("B3").Value = "GFTA_SIW_Complete"
$ws.Range("B4").Value = "GFTA_SIW_Incomplete_Reason"
$ws.Range("B5").Value = "GFTA_SIW_Valid"
$ws.Range("B6").Value = "GFTA_SIW_Invalid_Reason"

# ---------------------------------------------------------------------
# 3. Rows 7-9 used to be the "part 2 completed / reason / valid" question
#    block; that block now lives further down (rows 10-12), and rows 7-9
#    become the sounds-in-words Raw/Standard/Percentile rows (previously
#    rows 13-15, now deleted/replaced in place).
# ---------------------------------------------------------------------
$ws.Range("A7").Value = "Sounds-in-words Raw Score"
$ws.Range("B7").Value = "GFTA_SIW_Raw"
$ws.Range("C7").Value = "Numeric"
$ws.Range("D7").ClearContents()
$ws.Range("E7").ClearContents()

$ws.Range("A8").Value = "Sounds-in-words Standard Score"
$ws.Range("B8").Value = "GFTA_SIW_Stnd"
$ws.Range("C8").Value = "Numeric"
$ws.Range("D8").ClearContents()
$ws.Range("E8").ClearContents()

$ws.Range("A9").Value = "Sounds-in-words Percentile Rank"
$ws.Range("B9").Value = "GFTA_SIW_Perc"
$ws.Range("C9").Value = "Numeric"
$ws.Range("D9").ClearContents()
$ws.Range("E9").ClearContents()

# ---------------------------------------------------------------------
# 4. Rows 10-14: the "part 2" (sounds-in-sentences) question block, now
#    with its own complete/incomplete-reason/valid/invalid-reason rows
#    (mirroring part 1) plus the language-background follow-up.
# ---------------------------------------------------------------------
$ws.Range("A10").Value = "Was GFTA part 2 (Sounds-in-sentences) completed?"
$ws.Range("B10").Value = "GFTA_SIS_Complete"
$ws.Range("C10").Value = "Numeric"
$ws.Range("D10").Value = "0-1"
$ws.Range("E10").Value = "0=No, 1=Yes"

$ws.Range("A11").Value = "If no, reason:"
$ws.Range("B11").Value = "GFTA_SIS_Incomplete_Reason"
$ws.Range("C11").Value = "Text"
$ws.Range("D11").ClearContents()
$ws.Range("E11").ClearContents()

$ws.Range("A12").Value = "Was GFTA part 2 (Sounds-in-sentences) administration valid?"
$ws.Range("B12").Value = "GFTA_SIS_Valid"
$ws.Range("C12").Value = "Numeric"
$ws.Range("D12").Value = "0-1"
$ws.Range("E12").Value = "0=No, 1=Yes"

$ws.Range("A13").Value = "If no, reason:"
$ws.Range("B13").Value = "GFTA_SIS_Invalid_Reason"
$ws.Range("C13").Value = "Text"
$ws.Range("D13").ClearContents()
$ws.Range("E13").ClearContents()

$ws.Range("A14").Value = "If no, then possibly due to language background"
$ws.Range("B14").Value = "GFTA_SIS_Valid_Possible"
$ws.Range("C14").Value = "Text"
$ws.Range("D14").ClearContents()
$ws.Range("E14").ClearContents()

# ---------------------------------------------------------------------
# 5. Rows 15-17 (sounds-in-sentences Raw/Standard/Percentile) already line
#    up after the row-11 delete; nothing else to move.
# ---------------------------------------------------------------------

# ---------------------------------------------------------------------
# 6. Selection moves to B10 (matches the author's final cursor position).
# ---------------------------------------------------------------------
$ws.Range("B10").Select()
